$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of price data (Packham's Triumph pear, Región de O'Higgins) was
# collected. Insert 3 rows above the current top of this variety/quality
# block (row 940) so the whole block shifts down by three rows
# (940:965 -> 943:968), then populate the 3 freshly inserted rows with the
# new week's figures.
$ws.Rows("940:942").Insert()

# Row 940: Packham's Triumph / Especial
$ws.Range("A940").Value = 8
$ws.Range("B940").Value = "Terminal La Palmera de La Serena"
$ws.Range("C940").Value = "Coquimbo"
$ws.Range("D940").Value = 45041
$ws.Range("E940").Value = 4
$ws.Range("F940").Value = "Fruta"
$ws.Range("G940").Value = 100104
$ws.Range("H940").Value = "Frutos de pepita"
$ws.Range("I940").Value = 100104005
$ws.Range("J940").Value = "Pera"
$ws.Range("K940").Value = "Packham's Triumph"
$ws.Range("L940").Value = "Especial"
$ws.Range("M940").Value = 16
$ws.Range("N940").Value = 250000
$ws.Range("O940").Value = 260000
$ws.Range("P940").Value = 255000
$ws.Range("Q940").Value = "$/bins (450 kilos)"
$ws.Range("R940").Value = "Región de O'Higgins"
$ws.Range("S940").Value = 567
$ws.Range("T940").Value = 450

# Row 941: Packham's Triumph / Primera
$ws.Range("A941").Value = 8
$ws.Range("B941").Value = "Terminal La Palmera de La Serena"
$ws.Range("C941").Value = "Coquimbo"
$ws.Range("D941").Value = 45041
$ws.Range("E941").Value = 4
$ws.Range("F941").Value = "Fruta"
$ws.Range("G941").Value = 100104
$ws.Range("H941").Value = "Frutos de pepita"
$ws.Range("I941").Value = 100104005
$ws.Range("J941").Value = "Pera"
$ws.Range("K941").Value = "Packham's Triumph"
$ws.Range("L941").Value = "Primera"
$ws.Range("M941").Value = 12
$ws.Range("N941").Value = 220000
$ws.Range("O941").Value = 230000
$ws.Range("P941").Value = 225000
$ws.Range("Q941").Value = "$/bins (450 kilos)"
$ws.Range("R941").Value = "Región de O'Higgins"
$ws.Range("S941").Value = 500
$ws.Range("T941").Value = 450

# Row 942: Packham's Triumph / Segunda
$ws.Range("A942").Value = 8
$ws.Range("B942").Value = "Terminal La Palmera de La Serena"
$ws.Range("C942").Value = "Coquimbo"
$ws.Range("D942").Value = 45041
$ws.Range("E942").Value = 4
$ws.Range("F942").Value = "Fruta"
$ws.Range("G942").Value = 100104
$ws.Range("H942").Value = "Frutos de pepita"
$ws.Range("I942").Value = 100104005
$ws.Range("J942").Value = "Pera"
$ws.Range("K942").Value = "Packham's Triumph"
$ws.Range("L942").Value = "Segunda"
$ws.Range("M942").Value = 10
$ws.Range("N942").Value = 200000
$ws.Range("O942").Value = 210000
$ws.Range("P942").Value = 205000
$ws.Range("Q942").Value = "$/bins (450 kilos)"
$ws.Range("R942").Value = "Región de O'Higgins"
$ws.Range("S942").Value = 456
$ws.Range("T942").Value = 450
